# 20200317 Add the different situation email
# Update row 44 with corrected values and append new rows 45 and 46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper pattern: numeric-looking values that must stay stored as TEXT
# (the sheet uses text cells for the price/volume columns) are written
# by flipping the cell to a text number format first, assigning the
# value, then resetting the style back to Normal so no stray formatting
# is left behind.
# ---------------------------------------------------------------------

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ----------------------- Row 44 (existing row, values revised) -----------------------

Set-TextValue $ws.Range("D44") "120.00000000"
Set-TextValue $ws.Range("E44") "122.54000000"
Set-TextValue $ws.Range("F44") "1237674.88642000"
Set-TextValue $ws.Range("H44") "159192895.55261970"
$ws.Range("I44").Value = 398146
Set-TextValue $ws.Range("J44") "628481.86635000"
Set-TextValue $ws.Range("K44") "80877487.98928390"

$ws.Range("M44").Value = 122.5399999999998
$ws.Range("N44").Value = 128.3000000000001
$ws.Range("O44").Value = 165.9957142857143
$ws.Range("P44").Value = 199.1746666666667
$ws.Range("Q44").Value = 228.4786666666668
$ws.Range("R44").Value = 122.54
$ws.Range("S44").Value = 126.7756675048045
$ws.Range("T44").Value = 179.8633457903535
$ws.Range("U44").Value = 204.8401443871479
$ws.Range("V44").Value = -24.97679859679442
$ws.Range("W44").Value = -13.75694642572816
$ws.Range("X44").Value = -11.21985217106626

# ----------------------- Row 45 (new row) -----------------------

$ws.Range("A44").Copy()
$ws.Range("A45:A46").PasteSpecial(-4122)

$ws.Range("A45").Value = 43
Set-TextValue $ws.Range("B45") "122.54000000"
Set-TextValue $ws.Range("C45") "133.50000000"
Set-TextValue $ws.Range("D45") "120.16000000"
Set-TextValue $ws.Range("E45") "123.78000000"
Set-TextValue $ws.Range("F45") "1312950.90137000"
$ws.Range("G45").Value = 1584316799999
Set-TextValue $ws.Range("H45") "163436235.26909020"
$ws.Range("I45").Value = 346579
Set-TextValue $ws.Range("J45") "669525.68240000"
Set-TextValue $ws.Range("K45") "83359166.35016240"
Set-TextValue $ws.Range("L45") "2020-03-15 08:00:00"

$ws.Range("M45").Value = 123.7799999999998
$ws.Range("N45").Value = 123.1600000000001
$ws.Range("O45").Value = 155.1885714285714
$ws.Range("P45").Value = 192.946
$ws.Range("Q45").Value = 223.0996666666668
$ws.Range("R45").Value = 123.78
$ws.Range("S45").Value = 124.7785558349348
$ws.Range("T45").Value = 171.229592868324
$ws.Range("U45").Value = 198.6254195632536
$ws.Range("V45").Value = -27.39582669492958
$ws.Range("W45").Value = -16.48487100190949
$ws.Range("X45").Value = -10.91095569302009

# ----------------------- Row 46 (new row) -----------------------

$ws.Range("A46").Value = 44
Set-TextValue $ws.Range("B46") "123.82000000"
Set-TextValue $ws.Range("C46") "124.33000000"
Set-TextValue $ws.Range("D46") "101.10000000"
Set-TextValue $ws.Range("E46") "111.18000000"
Set-TextValue $ws.Range("F46") "1778230.76646000"
$ws.Range("G46").Value = 1584403199999
Set-TextValue $ws.Range("H46") "196383322.33181800"
$ws.Range("I46").Value = 466021
Set-TextValue $ws.Range("J46") "864592.43860000"
Set-TextValue $ws.Range("K46") "95415716.67891000"
Set-TextValue $ws.Range("L46") "2020-03-16 08:00:00"

$ws.Range("M46").Value = 111.1799999999998
$ws.Range("N46").Value = 117.48
$ws.Range("O46").Value = 142.0985714285714
$ws.Range("P46").Value = 185.8373333333333
$ws.Range("Q46").Value = 217.9763333333334
$ws.Range("R46").Value = 111.18
$ws.Range("S46").Value = 115.7128519449783
$ws.Range("T46").Value = 161.9861699288043
$ws.Range("U46").Value = 191.9384938731846
$ws.Range("V46").Value = -29.95232394438028
$ws.Range("W46").Value = -19.17847891357495
$ws.Range("X46").Value = -10.77384503080534
